# Update "想去人数" (wanted-to-go count) figures that changed between
# the two data refreshes, on both the "展览" and "全部类型" sheets.

$wb = $excel.ActiveWorkbook

$sheetNames = @("展览", "全部类型")

foreach ($name in $sheetNames) {
    $ws = $wb.Worksheets.Item($name)
    $ws.Range("F2").Value = 4757
    $ws.Range("F3").Value = 141
    $ws.Range("F4").Value = 829
}
